# Update the canonical terminology URLs (MOS/NOS hostname) and refresh the
# "Date" metadata timestamp, matching the upstream IG regeneration commit
# "Modif url canonique termino df9498eb894642b7264f6d5c9a38a249f1b02b34".
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Metadata sheet: "Date" row (A8 label / B8 value) -> new generation timestamp.
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# Elements sheet, "Binding Value Set" column (Z): esante.gouv.fr -> mos.esante.gouv.fr canonical URLs.
$wsElements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R36-AutreDiplomeObtenu/FHIR/TRE-R36-AutreDiplomeObtenu?vs"

# Column Z (Binding Value Set) widens slightly to fit the new, longer URLs.
$wsElements.Columns.Item(26).ColumnWidth = 80.7
